$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab14")

# Fix mojibake text in the Regional Economic Communities note (A103)
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Update recalculated statistics for "Africa, Fragile States" (row 97)
$ws.Range("C97").Value = 4.5540752227489802
$ws.Range("D97").Value = 0.34843122615263999
$ws.Range("E97").Value = 0.68496586726261999
$ws.Range("F97").Value = 0.69022881640837996
$ws.Range("G97").Value = 0.04513061001825
$ws.Range("H97").Value = 0.70881008184872996

# Update recalculated statistics for "ROW, Fragile States" (row 98)
$ws.Range("C98").Value = 5.2005227009455401
$ws.Range("D98").Value = 0.34283769006530002
$ws.Range("E98").Value = 0.73093719780444999
$ws.Range("F98").Value = 0.81080003082751995
$ws.Range("G98").Value = 0.091898612910880001
$ws.Range("H98").Value = 0.7763326416413
